$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.813.97"
$ws.Range("E2").Value = "  +0.15%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.080.41"
$ws.Range("E3").Value = "  -0.91%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.33"
$ws.Range("E5").Value = "  -0.05%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.625"
$ws.Range("E6").Value = "  +0.26%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "58.91"
$ws.Range("E7").Value = "  +2.25%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  +1.28%  "
$ws.Range("E10").Value = "  +0.91%  "
$ws.Range("E11").Value = "  +1.55%  "
$ws.Range("E12").Value = "  +1.84%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.12"
$ws.Range("E13").Value = "  -0.75%  "
$ws.Range("E15").Value = "  +2.29%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.131.48"
$ws.Range("E16").Value = "  +1.53%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "37.715.70"
$ws.Range("E17").Value = "  -0.03%  "
$ws.Range("E18").Value = "  -1.46%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "71.56"
$ws.Range("E19").Value = "  +1.35%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0837"
$ws.Range("E20").Value = "  +1.69%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "228.39"
$ws.Range("E21").Value = "  +0.28%  "
$ws.Range("E22").Value = "  -0.10%  "
$ws.Range("E23").Value = "  +0.87%  "
$ws.Range("E24").Value = "  -2.74%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "170.72"
$ws.Range("E25").Value = "  +1.66%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.21"
$ws.Range("E26").Value = "  +2.82%  "
$ws.Range("E27").Value = "  -2.68%  "
$ws.Range("E28").Value = "  +0.35%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.52"
$ws.Range("E29").Value = "  +0.10%  "
$ws.Range("E30").Value = "  +2.05%  "
$ws.Range("E31").Value = "  +0.93%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0636"
$ws.Range("E32").Value = "  +1.36%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.71"
$ws.Range("E33").Value = "  +2.07%  "
$ws.Range("E34").Value = "  -3.46%  "
$ws.Range("E35").Value = "  +0.04%  "
$ws.Range("E36").Value = "  -0.98%  "
$ws.Range("E37").Value = "  -0.08%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.37"
$ws.Range("E38").Value = "  -1.26%  "
$ws.Range("E39").Value = "  -0.58%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "99.75"
$ws.Range("E40").Value = "  +2.53%  "
$ws.Range("E41").Value = "  -2.38%  "
$ws.Range("E42").Value = "  +0.58%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "16.72"
$ws.Range("E43").Value = "  +5.98%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.441.94"
$ws.Range("E44").Value = "  -0.99%  "
$ws.Range("E45").Value = "  -0.43%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.18"
$ws.Range("E46").Value = "  +1.64%  "
$ws.Range("E47").Value = "  +0.43%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.43"
$ws.Range("E48").Value = "  +2.28%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.98"
$ws.Range("E49").Value = "  -1.14%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.270.38"
$ws.Range("E50").Value = "  -0.95%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "46.69"
$ws.Range("E51").Value = "  +0.81%  "
